$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NO_REGISTER cell (N2)
$ws.Range("N2").Value = "UP032303000071"

# Update the PREPARATION cell (F2) which embeds the register number in a multi-line note
$ws.Range("F2").Value = "Username : 31816;" + [char]10 + "Password : bni1234;" + [char]10 + "Role : 09 - Penyelia Settlement;" + [char]10 + "No Register : UP032303000071"

# Move selection to G2 to match the updated view state
$ws.Range("G2").Select()
